$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")
$ws.Activate()

# Week 9 (column J, week of 10/25) results recorded for the "Wookie Mistakes" table (rows 3-10)
$ws.Range("J3").Value  = "DNP"
$ws.Range("J4").Value  = "W"
$ws.Range("J5").Value  = "NA"
$ws.Range("J6").Value  = "W"
$ws.Range("J7").Value  = "DNP"
$ws.Range("J8").Value  = "L"
$ws.Range("J9").Value  = "L"
$ws.Range("J10").Value = "DNP"

# Week 9 (column J) results recorded for the "Safety Dance" table (rows 15-22)
$ws.Range("J15").Value = "L"
$ws.Range("J16").Value = "DNP"
$ws.Range("J17").Value = "W"
$ws.Range("J18").Value = "W"
$ws.Range("J19").Value = "L"
$ws.Range("J20").Value = "DNP"
$ws.Range("J21").Value = "DNP"
$ws.Range("J22").Value = "W"

# Reflect where the editor scrolled to / left the selection afterwards
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K18").Select()
